$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Memory Usage (bytes)"

$ws.Range("C2").Value = 16.44110679626465
$ws.Range("C3").Value = 16.0670280456543
$ws.Range("C4").Value = 16.10803604125977
$ws.Range("C5").Value = 16.05081558227539
$ws.Range("C6").Value = 15.72895050048828
